$d = $word.ActiveDocument

# Remove the leading blank paragraph
$d.Paragraphs.First.Range.Delete()

# Merge the split runs of "$Test=Yes2$" into a single run
$d.Content.Find.Execute("`$Test=Yes2`$", $false, $false, $false, $false, $false, $true, 1, $false, "`$Test=Yes2`$", 2) | Out-Null

# Append the new content block at the end of the document
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter("`r`r% **Complex Nested Fractions**`r`$I_N = \frac{\frac{V_{RC}}{R_2}}{\frac{V_E}{R_1}}`$`r`$P_{out} = \frac{\frac{F \cdot d}{T}}{\frac{I_E}{t}}`$`r`r% **Integral Expressions**`r`$\int_0^\infty x^2 \, dx`$`r`$\int_a^b f(x) \, dx = F(b) - F(a)`$`r`$A = \int_0^T v(t) \, dt`$`r`r% **Summation with Exponents**`r`$\sum_{n=1}^{\infty} \frac{1}{n^2} = \frac{\pi^2}{6}`$`r`$E = \sum_{n=0}^{\infty} e^{-n}`$`r`r% **Exponential and Logarithmic Functions**`r`$y = e^x \quad \text{and} \quad z = \log(x)`$`r`r% **Nested Exponents**`r`$y = e^{x^2 + \frac{1}{x}}`$`r`$z = x^{x^{x}}`$`r")

Write-Output $d.Paragraphs.Count
